$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values are forced to Text format so numeric-looking strings
# (e.g. "0.9967") are preserved as text instead of being parsed as numbers,
# matching the original inlineStr cell type. ClearFormats() afterwards
# removes the temporary "@" number-format style so no extra style index
# is introduced on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.666.53'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.647.94'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9967'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9976'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.42'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3803'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.18'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3607'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.245'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08203'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9967'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.49'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.528'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.369'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001230'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.643.96'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.00'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.05%  '
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.746'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.58'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9974'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.57'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.648.76'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("E25").Value = '  +2.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.121'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.25'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.44'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.204'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.89%  '
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.828.72'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.755'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.091'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +7.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.58'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +6.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.051'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -9.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02799'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08817'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.087'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.07026'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.80'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +4.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7062'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.329'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.84'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6506'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.84%  '
$ws.Range("E46").Value = '  +0.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9979'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.978'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07978'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.89'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.36%  '
$ws.Range("E51").Value = '  -0.72%  '
